$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 2 and fill it with the "дефИс / дЕфис" pair.
$ws.Rows("2:2").Insert()
$ws.Range("A2:B2").ClearFormats()
$ws.Range("A2").Value2 = "дефИс"
$ws.Range("B2").Value2 = "дЕфис"

# Insert another new row 2 (pushing the previous one down to row 3)
# and fill it with the "лыжнЯ / лЫжня" pair.
$ws.Rows("2:2").Insert()
$ws.Range("A2:B2").ClearFormats()
$ws.Range("A2").Value2 = "лыжнЯ"
$ws.Range("B2").Value2 = "лЫжня"

# Match the final selection/active cell from the authored workbook.
$ws.Range("E5").Select() | Out-Null
